# "full published + doi"
#
# The "Attributable deaths ... antibiotic-resistant bacteria ..." paper
# (previously listed only as "submitted", in the placeholder row that used
# to be at row 99) is now fully published. Its complete bibliographic
# details (volume, issue, pages, publication date, DOI) replace the old
# placeholder, and the row takes its rightful place in the date-sorted
# "papers" table (sorted ascending by column K, "date").
#
# Net effect on the papers table (rows 97-103), without changing the
# number of rows in the sheet:
#   - row 97 keeps its "published" formatting, but gets the new values
#     for the "Attributable deaths..." paper (full info, DOI, date 2019-01-01)
#   - row 98 gets row 97's old formatting, and the values that used to sit
#     in row 97 (the cadmium CKD paper)
#   - row 99 gets the values that used to sit in row 98 (the food-groups /
#     CHD paper); formatting is unchanged since rows 98/99 already share it
#   - rows 100/101/103 are untouched
#   - row 102 gains a DOI value in column N
#
# Also: the "chapters" sheet rows 5-7 gain their page ranges (from/to) and
# their publication date is corrected from 2019-01-01 to 2018-12-21.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "papers" (Tabel1)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("papers")

# Row 98 must take on the formatting that row 97 currently has (row 97
# itself keeps its own formatting - only its values change).
$ws.Range("A97:R97").Copy()
$ws.Range("A98:R98").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 97: now the fully published "Attributable deaths..." paper ---
$ws.Range("A97").Value = "Attributable deaths and disability-adjusted life-years caused by infections with antibiotic-resistant bacteria in the EU and the European Economic Area in 2015: a population-level modelling analysis"
$ws.Range("B97").Value = "Cassini, Alessandro; Högberg, Liselotte Diaz; Plachouras, Diamantis; Quattrocchi, Annalisa; Hoxha, Ana; Simonsen, Gunnar Skov; Colomb-Cotinat, Mélanie; Kretzschmar, Mirjam E.; Devleesschauwer, Brecht; Cecchini, Michele; Ouakrim, Driss Ait; Oliveira, Tiago Cravo; Struelens, Marc J.; Suetens, Carl; Monnet, Dominique L.; the Burden of AMR collaborative group"
$ws.Range("C97").Value = "The Lancet Infectious Diseases"
$ws.Range("D97").Value = "Lancet Infect. Dis."
$ws.Range("F97").Value = 19
# "issue" is stored as the text "1" (not the number 1) elsewhere in this
# column, so round-trip it through a scratch cell formatted as Text -
# a plain Value assignment here would be auto-coerced to a number and
# would also disturb G97's existing (bordered) cell style.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "1"
$ws.Range("Z1").Copy()
$ws.Range("G97").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("Z1").Clear()
$ws.Range("I97").Value = 56
$ws.Range("J97").Value = 66
$ws.Range("K97").Value = 43466
$ws.Range("N97").Value = "10.1016/S1473-3099(18)30605-4"

# --- Row 98: the cadmium CKD paper that used to be in row 97 ---
$ws.Range("A98").Value = "Global burden of late-stage chronic kidney disease resulting from dietary exposure to cadmium, 2015"
$ws.Range("B98").Value = "Zang, Yu; Devleesschauwer, Brecht; Bolger, P Michael; Goodman, Emily; Gibb, Herman J"
$ws.Range("C98").Value = "Environmental Research"
$ws.Range("D98").Value = "Environ. Res."
$ws.Range("E98").Value = 2019
$ws.Range("F98").Value = 169
$ws.Range("I98").Value = 72
$ws.Range("J98").Value = 78
$ws.Range("K98").Value = 43497
$ws.Range("N98").Value = "10.1016/j.envres.2018.10.005"

# --- Row 99: the food-groups / CHD paper that used to be in row 98 ---
$ws.Range("A99").Value = "Food groups and risk of coronary heart disease, stroke and heart failure: a systematic review and dose-response meta-analysis of prospective studies"
$ws.Range("B99").Value = "Bechthold, Angela; Boeing, Heiner; Schwedhelm, Carolina; Hoffmann, Georg; Knüppel, Sven; Iqbal, Khalid; De Henauw, Stefaan; Michels, Nathalie; Devleesschauwer, Brecht; Schlesinger, Sabrina; Schwingshackl, Lukas"
$ws.Range("C99").Value = "Critical Reviews in Food Science and Nutrition"
$ws.Range("D99").Value = "Crit. Rev. Food Sci. Nutr."
$ws.Range("K99").Value = 43831
$ws.Range("N99").Value = "10.1080/10408398.2017.1392288"

# --- Row 102: the methylmercury paper now has its DOI filled in ---
$ws.Range("N102").Value = "10.1016/j.envres.2018.12.042"

# ---------------------------------------------------------------------
# Sheet "chapters"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("chapters")

$ws2.Range("F5").Value = 83
$ws2.Range("G5").Value = 106
$ws2.Range("H5").Value = 43455

$ws2.Range("F6").Value = 107
$ws2.Range("G6").Value = 122
$ws2.Range("H6").Value = 43455

$ws2.Range("F7").Value = 143
$ws2.Range("G7").Value = 174
$ws2.Range("H7").Value = 43455
